$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text so numeric-looking values are not
# auto-converted to numbers by Excel (keeps them as inline/shared strings).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.547.53"
$ws.Range("E2").Value = "  -0.81%  "

$ws.Range("D3").Value = "2.710.55"
$ws.Range("E3").Value = "  -1.91%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "558.49"
$ws.Range("E5").Value = "  -3.16%  "

$ws.Range("D6").Value = "156.78"
$ws.Range("E6").Value = "  -1.56%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  -2.53%  "

$ws.Range("D9").Value = "0.106"
$ws.Range("E9").Value = "  -3.29%  "

$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  +0.04%  "

$ws.Range("D11").Value = "5.51"
$ws.Range("E11").Value = "  -5.25%  "

$ws.Range("D12").Value = "0.371"
$ws.Range("E12").Value = "  -3.69%  "

$ws.Range("D13").Value = "3.193.61"
$ws.Range("E13").Value = "  -1.69%  "

$ws.Range("D14").Value = "26.48"
$ws.Range("E14").Value = "  -1.69%  "

$ws.Range("D15").Value = "63.391.86"
$ws.Range("E15").Value = "  -0.45%  "

$ws.Range("D16").Value = "0.0000145"
$ws.Range("E16").Value = "  -3.64%  "

$ws.Range("D17").Value = "2.719.32"
$ws.Range("E17").Value = "  -1.58%  "

$ws.Range("D18").Value = "12.15"
$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("D19").Value = "4.62"
$ws.Range("E19").Value = "  -4.61%  "

$ws.Range("D20").Value = "349.15"
$ws.Range("E20").Value = "  -1.97%  "

$ws.Range("D21").Value = "6.40"
$ws.Range("E21").Value = "  -3.94%  "

$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.21%  "

$ws.Range("D23").Value = "0.509"
$ws.Range("E23").Value = "  -3.39%  "

$ws.Range("D24").Value = "64.06"
$ws.Range("E24").Value = "  -1.52%  "

$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("D27").Value = "8.14"
$ws.Range("E27").Value = "  -4.59%  "

$ws.Range("D28").Value = "0.0₃0870"
$ws.Range("E28").Value = "  -4.10%  "

$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "1.93"
$ws.Range("E29").Value = "  -0.96%  "

$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "1.34"
$ws.Range("E30").Value = "  +5.90%  "

$ws.Range("D31").Value = "7.12"
$ws.Range("E31").Value = "  -2.81%  "

$ws.Range("D32").Value = "165.34"
$ws.Range("E32").Value = "  -2.62%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "19.74"
$ws.Range("E34").Value = "  -2.06%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "4.81"
$ws.Range("E35").Value = "  -2.34%  "

$ws.Range("D36").Value = "1.44"
$ws.Range("E36").Value = "  -1.71%  "

$ws.Range("D37").Value = "1.78"
$ws.Range("E37").Value = "  -1.40%  "

$ws.Range("D38").Value = "344.36"
$ws.Range("E38").Value = "  -1.70%  "

$ws.Range("D39").Value = "0.951"
$ws.Range("E39").Value = "  -5.42%  "

$ws.Range("D40").Value = "6.09"
$ws.Range("E40").Value = "  -2.85%  "

$ws.Range("D41").Value = "3.97"
$ws.Range("E41").Value = "  -4.59%  "

$ws.Range("D42").Value = "38.32"
$ws.Range("E42").Value = "  -1.99%  "

$ws.Range("D43").Value = "21.01"
$ws.Range("E43").Value = "  -3.72%  "

$ws.Range("D44").Value = "20.51"
$ws.Range("E44").Value = "  -4.40%  "

$ws.Range("D45").Value = "0.626"
$ws.Range("E45").Value = "  -0.65%  "

$ws.Range("D46").Value = "0.0566"
$ws.Range("E46").Value = "  -3.43%  "

$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  +0.12%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "132.14"
$ws.Range("E48").Value = "  -4.00%  "

$ws.Range("E49").Value = "  +0.43%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.0978"
$ws.Range("E50").Value = "  -3.28%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.120.44"
$ws.Range("E51").Value = "  -0.07%  "

# Revert the temporary text formatting on column D so the saved styles
# match the original (no explicit number-format/style on these cells).
$ws.Range("D2:D51").ClearFormats()
